$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.864.67"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.99%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.032.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.14%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.84"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.81"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.80%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.027.12"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.00%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.44"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +8.23%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.466"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000234"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.49"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.60%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.533.88"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.10"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.848.01"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.030.29"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "451.49"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.30"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.695"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.03%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.51%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.30"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.03"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +5.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.32"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.77%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.49"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +6.79%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.90%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.55"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.12%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0866"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.51%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.92"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.18"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +10.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.11"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.50"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.31%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.09"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.34%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +14.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.63"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "394.54"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.01%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.729.61"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.31"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.35"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.63%  "
